$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new data adds two price-report rows for "Bruselas (repollito)" at the
# Vega Modelo de Temuco market. One row is inserted before the existing
# row 82 (date 2022-06-06 / serial 44699), the other is inserted before the
# existing row 90 (date 2021-08-13 / serial 44421) - both pushing the
# following rows down.

# --- Insert first new row at row 82 ---
$ws.Rows.Item(82).Insert()

$ws.Range("A82").Value = 10
$ws.Range("B82").Value = "Vega Modelo de Temuco"
$ws.Range("C82").Value = "La Araucanía"
$ws.Range("D82").Value = 44748
$ws.Range("E82").Value = 9
$ws.Range("F82").Value = 100112035
$ws.Range("G82").Value = "Bruselas (repollito)"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 45
$ws.Range("K82").Value = 26000
$ws.Range("L82").Value = 26000
$ws.Range("M82").Value = 26000
$ws.Range("N82").Value = "$/malla 10 kilos"
$ws.Range("O82").Value = "Provincia de Quillota"
$ws.Range("P82").Value = 2600
$ws.Range("Q82").Value = 10
$ws.Range("R82").Value = "Hortaliza"

# --- Insert second new row at (what is now) row 90 ---
$ws.Rows.Item(90).Insert()

$ws.Range("A90").Value = 10
$ws.Range("B90").Value = "Vega Modelo de Temuco"
$ws.Range("C90").Value = "La Araucanía"
$ws.Range("D90").Value = 44747
$ws.Range("E90").Value = 9
$ws.Range("F90").Value = 100112035
$ws.Range("G90").Value = "Bruselas (repollito)"
$ws.Range("H90").Value = "Sin especificar"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 20
$ws.Range("K90").Value = 28000
$ws.Range("L90").Value = 28000
$ws.Range("M90").Value = 28000
$ws.Range("N90").Value = "$/malla 10 kilos"
$ws.Range("O90").Value = "Provincia de Quillota"
$ws.Range("P90").Value = 2800
$ws.Range("Q90").Value = 10
$ws.Range("R90").Value = "Hortaliza"
